# Commit: table style GUID change on 3 tables + slide-master theme recolor
# (Integral/"Red Violet" colours -> stock "Office" colours); theme1.xml
# keeps its part name/relationships, only the 12 theme colour slots change.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-apply the table style on every table in the deck (slides 14, 15, 16 -
#    shape 1 in each case) from the custom "Table_0" style to the new style
#    id. Table.Style is read-only as a plain property assignment in this
#    object model; ApplyStyle() is the supported mutator.
# ---------------------------------------------------------------------------
$newTableStyleId = "{62C8B24B-CFEE-412F-AF64-E08759011711}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Recolour the presentation's theme (ppt/theme/theme1.xml, used by
#    slideMaster1 / all slides) from the "Integral" / "Red Violet" palette
#    to the stock Office theme palette. ThemeColorScheme indices 1-12 map
#    onto dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink in that order.
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

# RGB() in the PowerPoint object model packs as 0x00BBGGRR, so convert from
# the usual 0xRRGGBB hex triplet before assigning.
function ToVbRgb([int]$rrggbb) {
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le 12; $k++) {
    $themeColorScheme.Item($k).RGB = ToVbRgb($officeThemeColors[$k - 1])
}
